$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes for columns M (13) and N (14): 25 -> 24
$ws.Columns.Item(13).ColumnWidth = 23.166666666666668
$ws.Columns.Item(14).ColumnWidth = 23.166666666666668

# Cell value updates (sensor data refresh + banner timestamps)
$ws.Range('E2').Value = '2026-02-06 01:17:45'
$ws.Range('M2').Value = '-0.4 °C 0:54 TU'
$ws.Range('N2').Value = '-0.7 °C 0:41 TU'
$ws.Range('E3').Value = '2026-02-06 01:17:48'
$ws.Range('H3').Value = '87%'
$ws.Range('L3').Value = '39.2 km/h - 262º 0:54 TU'
$ws.Range('M3').Value = '-1.5 °C 0:57 TU'
$ws.Range('O3').Value = '-2.4 °C'
$ws.Range('E4').Value = '2026-02-06 01:17:50'
$ws.Range('H4').Value = '51%'
$ws.Range('J4').Value = '990.9 hPa'
$ws.Range('N4').Value = '14.0 °C 0:55 TU'
$ws.Range('O4').Value = '14.6 °C'
$ws.Range('E5').Value = '2026-02-06 01:17:52'
$ws.Range('J5').Value = '991.6 hPa'
$ws.Range('L5').Value = '13.0 km/h - 340º 0:59 TU'
$ws.Range('E6').Value = '2026-02-06 01:17:55'
$ws.Range('H6').Value = '50%'
$ws.Range('J6').Value = '993.1 hPa'
$ws.Range('L6').Value = '37.8 km/h - 315º 0:59 TU'
$ws.Range('M6').Value = '15.4 °C 0:57 TU'
$ws.Range('O6').Value = '15.0 °C'
$ws.Range('E7').Value = '2026-02-06 01:17:58'
$ws.Range('J7').Value = '993.0 hPa'
$ws.Range('L7').Value = '36.4 km/h - 247º 0:58 TU'
$ws.Range('E8').Value = '2026-02-06 01:18:00'
$ws.Range('N8').Value = '7.2 °C 0:38 TU'
$ws.Range('E9').Value = '2026-02-06 01:18:03'
$ws.Range('E10').Value = '2026-02-06 01:18:05'
$ws.Range('M10').Value = '6.1 °C 0:34 TU'
$ws.Range('O10').Value = '5.7 °C'
$ws.Range('E11').Value = '2026-02-06 01:18:07'
$ws.Range('H11').Value = '85%'
$ws.Range('J11').Value = '994.3 hPa'
$ws.Range('N11').Value = '4.7 °C 0:32 TU'
$ws.Range('O11').Value = '4.9 °C'
$ws.Range('E12').Value = '2026-02-06 01:18:10'
$ws.Range('H12').Value = '60%'
$ws.Range('N12').Value = '12.2 °C 0:40 TU'
$ws.Range('O12').Value = '13.2 °C'
$ws.Range('E13').Value = '2026-02-06 01:18:12'
$ws.Range('H13').Value = '84%'
$ws.Range('M13').Value = '9.3 °C 0:53 TU'
$ws.Range('O13').Value = '8.2 °C'
$ws.Range('E14').Value = '2026-02-06 01:18:15'
$ws.Range('H14').Value = '76%'
$ws.Range('I14').Value = '0.1 mm'
$ws.Range('L14').Value = '66.2 km/h - 200º 0:40 TU'
$ws.Range('M14').Value = '-3.1 °C 0:54 TU'
$ws.Range('E15').Value = '2026-02-06 01:18:17'
$ws.Range('H15').Value = '63%'
$ws.Range('J15').Value = '991.5 hPa'
$ws.Range('N15').Value = '8.7 °C 0:57 TU'
$ws.Range('O15').Value = '11.8 °C'
$ws.Range('E16').Value = '2026-02-06 01:18:20'
$ws.Range('N16').Value = '4.0 °C 0:51 TU'
$ws.Range('O16').Value = '4.1 °C'
$ws.Range('E17').Value = '2026-02-06 01:18:22'
$ws.Range('M17').Value = '3.9 °C 0:55 TU'
$ws.Range('N17').Value = '3.1 °C 0:30 TU'
$ws.Range('O17').Value = '3.4 °C'
$ws.Range('E18').Value = '2026-02-06 01:18:25'
$ws.Range('G18').Value = '119 cm'
$ws.Range('L18').Value = '26.3 km/h - 300º 0:53 TU'
$ws.Range('N18').Value = '-4.4 °C 0:40 TU'
$ws.Range('O18').Value = '-4.3 °C'
$ws.Range('E19').Value = '2026-02-06 01:18:27'
$ws.Range('H19').Value = '97%'
$ws.Range('J19').Value = '996.2 hPa'
$ws.Range('L19').Value = '19.4 km/h - 289º 0:39 TU'
$ws.Range('M19').Value = '8.5 °C 0:50 TU'
$ws.Range('O19').Value = '7.4 °C'
$ws.Range('E20').Value = '2026-02-06 01:18:30'
$ws.Range('H20').Value = '74%'
$ws.Range('L20').Value = '28.8 km/h - 260º 0:35 TU'
$ws.Range('M20').Value = '-1.0 °C 0:44 TU'
$ws.Range('N20').Value = '-2.2 °C 0:56 TU'
$ws.Range('O20').Value = '-1.6 °C'
$ws.Range('E21').Value = '2026-02-06 01:18:32'
$ws.Range('H21').Value = '73%'
$ws.Range('I21').Value = '0.0 mm'
$ws.Range('J21').Value = '992.0 hPa'
$ws.Range('K21').Value = '0.0 MJ/m2'
$ws.Range('L21').Value = '14.4 km/h - 178º 0:16 TU'
$ws.Range('M21').Value = '9.7 °C 0:12 TU'
$ws.Range('N21').Value = '5.8 °C 0:58 TU'
$ws.Range('O21').Value = '7.9 °C'
$ws.Range('E22').Value = '2026-02-06 01:18:34'
$ws.Range('E23').Value = '2026-02-06 01:18:37'
$ws.Range('J23').Value = '992.1 hPa'
$ws.Range('L23').Value = '14.4 km/h - 44º 0:43 TU'
$ws.Range('N23').Value = '6.8 °C 0:56 TU'
$ws.Range('O23').Value = '7.3 °C'
$ws.Range('E24').Value = '2026-02-06 01:18:39'
$ws.Range('H24').Value = '64%'
$ws.Range('J24').Value = '991.4 hPa'
$ws.Range('L24').Value = '19.1 km/h - 10º 0:57 TU'
$ws.Range('N24').Value = '10.7 °C 0:57 TU'
$ws.Range('O24').Value = '12.1 °C'
$ws.Range('E25').Value = '2026-02-06 01:18:42'
$ws.Range('H25').Value = '91%'
$ws.Range('J25').Value = '994.7 hPa'
$ws.Range('N25').Value = '2.0 °C 0:33 TU'
$ws.Range('O25').Value = '2.5 °C'
$ws.Range('E26').Value = '2026-02-06 01:18:44'
$ws.Range('H26').Value = '84%'
$ws.Range('L26').Value = '20.5 km/h - 21º 0:34 TU'
$ws.Range('N26').Value = '-0.2 °C 0:49 TU'
$ws.Range('E27').Value = '2026-02-06 01:18:47'
$ws.Range('H27').Value = '92%'
$ws.Range('J27').Value = '991.7 hPa'
$ws.Range('N27').Value = '6.8 °C 0:58 TU'
$ws.Range('O27').Value = '8.9 °C'
$ws.Range('E28').Value = '2026-02-06 01:18:49'
$ws.Range('J28').Value = '993.3 hPa'
$ws.Range('E29').Value = '2026-02-06 01:18:51'
$ws.Range('L29').Value = '55.8 km/h - 264º 0:37 TU'
$ws.Range('E30').Value = '2026-02-06 01:18:54'
$ws.Range('H30').Value = '67%'
$ws.Range('N30').Value = '-3.0 °C 0:55 TU'
$ws.Range('O30').Value = '-2.2 °C'
$ws.Range('E31').Value = '2026-02-06 01:18:56'
$ws.Range('J31').Value = '995.9 hPa'
$ws.Range('N31').Value = '4.8 °C 0:59 TU'
$ws.Range('E32').Value = '2026-02-06 01:18:59'
$ws.Range('J32').Value = '993.7 hPa'
$ws.Range('L32').Value = '65.2 km/h - 295º 0:56 TU'
$ws.Range('E33').Value = '2026-02-06 01:19:01'
$ws.Range('H33').Value = '93%'
$ws.Range('N33').Value = '7.5 °C 0:30 TU'
$ws.Range('O33').Value = '8.3 °C'
$ws.Range('E34').Value = '2026-02-06 01:19:03'
$ws.Range('H34').Value = '71%'
$ws.Range('M34').Value = '10.1 °C 0:51 TU'
$ws.Range('O34').Value = '9.9 °C'
$ws.Range('E35').Value = '2026-02-06 01:19:06'
$ws.Range('N35').Value = '-2.8 °C 0:49 TU'
$ws.Range('E36').Value = '2026-02-06 01:19:08'
$ws.Range('J36').Value = '994.9 hPa'
